$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.326.04'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.412.28'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.07'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.42'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  +4.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.422.48'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.19'
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.999.80'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  -4.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.82'
$ws.Range('E16').Value = '  -4.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.389.51'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.393.64'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.93'
$ws.Range('E20').Value = '  -4.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '373.89'
$ws.Range('E21').Value = '  -4.74%  '
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.17'
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000118'
$ws.Range('E26').Value = '  -5.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.04'
$ws.Range('E27').Value = '  +5.15%  '
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('E30').Value = '  +2.33%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.06'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.20'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('E35').Value = '  +5.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.47'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0756'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.78'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.66'
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.847.02'
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.61'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.68'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.09'
$ws.Range('E44').Value = '  +7.73%  '
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.767'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '315.91'
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.07'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.108'
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.853'
$ws.Range('E51').Value = '  -2.67%  '
